$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("JSON_Data")

# Read source values from JSON_Data before inserting the new sheet
$values = @()
for ($r = 1; $r -le 6; $r++) {
    $values += , $ws1.Cells.Item($r, 1).Value2
}

# Add the new "DB_Data" sheet right after JSON_Data
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "DB_Data"

# Populate DB_Data: rows 1-6 mirror JSON_Data, row 7 gets the new e-mail value
for ($r = 1; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 1).Value2 = $values[$r - 1]
}
$ws2.Cells.Item(7, 1).Value2 = "pkailasam@deloitte"

# Highlight JSON_Data!A7 with a red fill
$ws1.Range("A7").Interior.ColorIndex = 3

# Keep JSON_Data as the active/selected sheet (as it was before the edit)
$ws1.Activate()
